$wb = $excel.ActiveWorkbook

# ALC row 17 (G=38956)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1249.8077
$ws.Range("I17").Value = 670
$ws.Range("J17").Value = 1355.2273
$ws.Range("K17").Value = 2010
$ws.Range("L17").Value = 4065.6819
$ws.Range("M17").Value = -1842
$ws.Range("N17").Value = -4401.6819

# ALC row 18 (G=5471)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 142869730
$ws.Range("I18").Value = 333335360
$ws.Range("K18").Value = 333335360
$ws.Range("M18").Value = -333335076

# ALC row 51 (G=5486)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 29999.5
$ws.Range("J51").Value = 29999.5
$ws.Range("L51").Value = 29999.5
$ws.Range("N51").Value = -30967.5

# ALC row 62 (G=27781)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 23718.818
$ws.Range("I62").Value = 14974.875
$ws.Range("K62").Value = 14974.875
$ws.Range("M62").Value = -14350.875

# ALC row 65 (G=27781)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 23718.818
$ws.Range("I65").Value = 14974.875
$ws.Range("K65").Value = 74874.375
$ws.Range("M65").Value = -71754.375

# ALC row 88 (G=12608)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 63016910
$ws.Range("I88").Value = 166667540
$ws.Range("K88").Value = 166667540
$ws.Range("M88").Value = -166667134

# ALC row 91 (G=12608)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 63016910
$ws.Range("I91").Value = 166667540
$ws.Range("K91").Value = 166667540
$ws.Range("M91").Value = -166666136

# ALC row 132 (G=44049)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1164.25
$ws.Range("I132").Value = 1139.6154
$ws.Range("K132").Value = 3418.8462
$ws.Range("M132").Value = -888.8462

# ALC row 137 (G=44013)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4353.243
$ws.Range("J137").Value = 6496.7144
$ws.Range("L137").Value = 19490.1432
$ws.Range("N137").Value = -24590.1432

# ARM row 2 (G=27713)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 62502308
$ws.Range("I2").Value = 2026.909
$ws.Range("J2").Value = 200002930
$ws.Range("K2").Value = 2026.909
$ws.Range("L2").Value = 200002930
$ws.Range("M2").Value = -1913.909
$ws.Range("N2").Value = -200003156

# ARM row 7 (G=27125)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("N7").Value = 0

# ARM row 32 (G=44147)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3338168.2
$ws.Range("I32").Value = 3574645
$ws.Range("J32").Value = 27495.5
$ws.Range("K32").Value = 3574645
$ws.Range("L32").Value = 27495.5
$ws.Range("M32").Value = -3574358
$ws.Range("N32").Value = -28069.5

# ARM row 61 (G=43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 10989.154
$ws.Range("I61").Value = 1711
$ws.Range("K61").Value = 1711
$ws.Range("M61").Value = -1499

# ARM row 74 (G=44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 87592.75
$ws.Range("I74").Value = 251649.75
$ws.Range("J74").Value = 5564.25
$ws.Range("K74").Value = 251649.75
$ws.Range("L74").Value = 5564.25
$ws.Range("M74").Value = -250775.75
$ws.Range("N74").Value = -7312.25

# ARM row 77 (G=44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 87592.75
$ws.Range("I77").Value = 251649.75
$ws.Range("J77").Value = 5564.25
$ws.Range("K77").Value = 1258248.75
$ws.Range("L77").Value = 27821.25
$ws.Range("M77").Value = -1253880.75
$ws.Range("N77").Value = -36557.25

# ARM row 97 (G=19941)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2526070
$ws.Range("I97").Value = 654.8214
$ws.Range("J97").Value = 16668395
$ws.Range("K97").Value = 654.8214
$ws.Range("L97").Value = 16668395
$ws.Range("M97").Value = -158.8214
$ws.Range("N97").Value = -16669387

# ARM row 116 (G=27713)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 62502308
$ws.Range("I116").Value = 2026.909
$ws.Range("J116").Value = 200002930
$ws.Range("K116").Value = 2026.909
$ws.Range("L116").Value = 200002930
$ws.Range("M116").Value = 267.0909999999999
$ws.Range("N116").Value = -200007518

# ARM row 136 (G=43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 10989.154
$ws.Range("I136").Value = 1711
$ws.Range("K136").Value = 5133
$ws.Range("M136").Value = -2583

# ARM row 140 (G=42496)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H140").Value = 49999.5
$ws.Range("J140").Value = 49999.5
$ws.Range("N140").Value = -60359.5

# BSM row 3 (G=27713)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 62502308
$ws.Range("I3").Value = 2026.909
$ws.Range("J3").Value = 200002930
$ws.Range("K3").Value = 2026.909
$ws.Range("L3").Value = 200002930
$ws.Range("M3").Value = -1912.909
$ws.Range("N3").Value = -200003158

# BSM row 13 (G=27127)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 70000
$ws.Range("J13").Value = 70000
$ws.Range("L13").Value = 70000
$ws.Range("N13").Value = -70336

# BSM row 60 (G=43232)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 48351.832
$ws.Range("J60").Value = 48351.832
$ws.Range("L60").Value = 48351.832
$ws.Range("N60").Value = -49549.832

# BSM row 74 (G=43246)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 66853.336
$ws.Range("J74").Value = 66853.336
$ws.Range("L74").Value = 66853.336
$ws.Range("N74").Value = -68725.336

# BSM row 77 (G=43246)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H77").Value = 66853.336
$ws.Range("J77").Value = 66853.336
$ws.Range("L77").Value = 200560.008
$ws.Range("N77").Value = -209920.008

# BSM row 86 (G=12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 101226.2
$ws.Range("I86").Value = 125926.625
$ws.Range("K86").Value = 125926.625
$ws.Range("M86").Value = -124803.625

# BSM row 89 (G=12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 101226.2
$ws.Range("I89").Value = 125926.625
$ws.Range("K89").Value = 629633.125
$ws.Range("M89").Value = -624017.125

# BSM row 134 (G=43998)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6679.242
$ws.Range("I134").Value = 2284.7144
$ws.Range("J134").Value = 9917.315000000001
$ws.Range("K134").Value = 6854.1432
$ws.Range("L134").Value = 29751.945
$ws.Range("M134").Value = -4319.1432
$ws.Range("N134").Value = -34821.945

# CRP row 16 (G=27691)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5528.375
$ws.Range("I16").Value = 1824
$ws.Range("K16").Value = 1824
$ws.Range("M16").Value = -1537

# CRP row 36 (G=1845)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 26198.75
$ws.Range("I36").Value = 8998.333000000001
$ws.Range("J36").Value = 31932.223
$ws.Range("K36").Value = 8998.333000000001
$ws.Range("L36").Value = 31932.223
$ws.Range("M36").Value = -8610.333000000001
$ws.Range("N36").Value = -32708.223

# CRP row 40 (G=1845)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H40").Value = 26198.75
$ws.Range("I40").Value = 8998.333000000001
$ws.Range("J40").Value = 31932.223
$ws.Range("K40").Value = 8998.333000000001
$ws.Range("L40").Value = 31932.223
$ws.Range("M40").Value = -8838.333000000001
$ws.Range("N40").Value = -32252.223

# CRP row 113 (G=27691)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 5528.375
$ws.Range("I113").Value = 1824
$ws.Range("K113").Value = 1824
$ws.Range("M113").Value = 346

# CRP row 122 (G=36196)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 4696.375
$ws.Range("I122").Value = 3458.3
$ws.Range("K122").Value = 10374.9
$ws.Range("M122").Value = -7924.900000000001

# CRP row 134 (G=44020)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 6783.057
$ws.Range("I134").Value = 3077.7
$ws.Range("K134").Value = 9233.099999999999
$ws.Range("M134").Value = -6698.099999999999

# CUL row 41 (G=4700)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 800
$ws.Range("I41").Value = 800
$ws.Range("K41").Value = 2400
$ws.Range("M41").Value = -2062

# CUL row 97 (G=19846)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 548.125
$ws.Range("I97").Value = 394
$ws.Range("K97").Value = 1182
$ws.Range("M97").Value = -686

# CUL row 131 (G=36060)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2240.5925
$ws.Range("I131").Value = 1269.1538
$ws.Range("J131").Value = 2548.6099
$ws.Range("K131").Value = 3807.4614
$ws.Range("L131").Value = 7645.8297
$ws.Range("M131").Value = 1232.5386
$ws.Range("N131").Value = -17725.8297

# GSM row 102 (G=36169)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3386.6667
$ws.Range("I102").Value = 3509.318
$ws.Range("K102").Value = 3509.318
$ws.Range("M102").Value = -1887.318

# GSM row 106 (G=18722)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H106").Value = 53884.5
$ws.Range("J106").Value = 53884.5
$ws.Range("N106").Value = -56408.5

# LTW row 7 (G=36249)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6541.3335
$ws.Range("I7").Value = 4218
$ws.Range("J7").Value = 8400
$ws.Range("K7").Value = 4218
$ws.Range("L7").Value = 8400
$ws.Range("M7").Value = -4106
$ws.Range("N7").Value = -8624

# LTW row 46 (G=5282)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4633137.5
$ws.Range("I46").Value = 2582.6
$ws.Range("K46").Value = 2582.6
$ws.Range("M46").Value = -2394.6

# LTW row 55 (G=5284)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 333.65384
$ws.Range("I55").Value = 18.583334
$ws.Range("J55").Value = 603.7143
$ws.Range("K55").Value = 18.583334
$ws.Range("L55").Value = 603.7143
$ws.Range("M55").Value = 154.416666
$ws.Range("N55").Value = -949.7143

# LTW row 126 (G=36249)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 6541.3335
$ws.Range("I126").Value = 4218
$ws.Range("J126").Value = 8400
$ws.Range("K126").Value = 12654
$ws.Range("L126").Value = 25200
$ws.Range("M126").Value = -10184
$ws.Range("N126").Value = -30140

# WVR row 41 (G=21725)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 26806.2
$ws.Range("J41").Value = 26686.5
$ws.Range("L41").Value = 26686.5
$ws.Range("N41").Value = -27466.5

# WVR row 42 (G=3372)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").Value = ""

# WVR row 122 (G=36208)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 15753974
$ws.Range("I122").Value = 21915584
$ws.Range("J122").Value = 7634.1113
$ws.Range("K122").Value = 65746752
$ws.Range("L122").Value = 22902.3339
$ws.Range("M122").Value = -65744302
$ws.Range("N122").Value = -27802.3339
